$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the stray _GoBack bookmark that sits at the end of the
#    "Machines adjustments" bullet in the body text.
# ------------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# ------------------------------------------------------------------
# 2. Rework the header. Touching the section's headers/footers makes
#    Word materialise the full even/default/first header & footer
#    parts (header1-3.xml, footer1-3.xml) and wires up the sectPr
#    headerReference/footerReference set automatically.
#
#    Index mapping for Section.Headers / Section.Footers:
#      1 = wdHeaderFooterPrimary (default)  -> header2.xml / footer2.xml
#      2 = wdHeaderFooterFirstPage (first)  -> header3.xml / footer3.xml
#      3 = wdHeaderFooterEvenPages (even)   -> header1.xml / footer1.xml
# ------------------------------------------------------------------
$sec = $d.Sections.Item(1)

# Touch the even header/footer so the even/first parts are minted too
# (they stay empty, matching the target document).
$evenHeader = $sec.Headers.Item(3)
$evenFooter = $sec.Footers.Item(3)
$firstHeader = $sec.Headers.Item(2)
$firstFooter = $sec.Footers.Item(2)
$defaultFooter = $sec.Footers.Item(1)

# Primary (default) header: update the title text and re-home the
# _GoBack bookmark into its own trailing paragraph, exactly as Word
# leaves it after the last edit was made inside the header.
$hdrPrimary = $sec.Headers.Item(1)
$hdrPrimary.Range.Text = "L6 – Group 8: Meeting minutes`rXMARKER"

$r = $hdrPrimary.Range
$r.Find.Execute("XMARKER", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(1)
$r.Bookmarks.Add("_GoBack")

$r2 = $hdrPrimary.Range
$r2.Find.Execute("XMARKER", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)
